$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("NCTId" shifts from C to D, etc.)
# This mirrors Excel's "Insert Column" behaviour: every column from C
# onward moves one to the right, and the new column C inherits the
# formatting of the old column C (so the header keeps its bold style).
$ws.Columns("C:C").Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 3).Value = "statut_name"

# Fill in the new "statut_name" column for every data row using the
# (now shifted) results_1y / results_3y / results columns (J, K, L).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $results1y = $ws.Cells.Item($r, 10).Value2
    $results3y = $ws.Cells.Item($r, 11).Value2
    $results   = $ws.Cells.Item($r, 12).Value2

    if ($results1y -eq $true) {
        $statutName = "résultat et / ou publication posté dans les 12 mois"
    } elseif ($results3y -eq $true) {
        $statutName = "résultat et / ou publication posté dans les 36 mois"
    } elseif ($results -eq $true) {
        $statutName = "résultat et / ou publication posté"
    } else {
        $statutName = "pas de résultat ni de publication"
    }

    $ws.Cells.Item($r, 3).Value = $statutName
}
